# Updated initial peak filtering algorithm - row 48 (TMAO-d9 related entry) values
# were recalculated with the new stricter rmt-window iteration logic, and the
# RMT tolerance percent cell (G48) now carries the same number format as the
# rest of that column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# G48 should use the same style/number-format ("0.00", style index 6) as every
# other cell in column G instead of the plain general style it had before.
$ws.Range("G48").NumberFormat = "0.00"

# Recalculated peak.N.rmt values for row 48 (columns I:U) from the reapplied,
# more strict rmt-window peak filtering pass.
$ws.Range("I48").Value = 1.2735210561497325
$ws.Range("J48").Value = 1.2600774403912778
$ws.Range("K48").Value = 1.2488874761602033
$ws.Range("L48").Value = 1.2387480097860277
$ws.Range("M48").Value = 1.2301240939622784
$ws.Range("N48").Value = 1.2175800526882119
$ws.Range("O48").Value = 1.209126854572496
$ws.Range("P48").Value = 1.2000355239786857
$ws.Range("Q48").Value = 1.1905308165606423
$ws.Range("R48").Value = 1.1828998056796101
$ws.Range("S48").Value = 1.1759927797833933
$ws.Range("T48").Value = 1.1682678120594119
$ws.Range("U48").Value = 1.1612592044760282

# Reflect where the author had scrolled to / selected when saving: the view
# had scrolled down so row 34 is the top visible row, and the selection spans
# F52:H67 with F52 as the active cell.
$excel.ActiveWindow.ScrollColumn = 1
$excel.ActiveWindow.ScrollRow = 34
$ws.Range("F52:H67").Select()
